$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Column D (MyForecast)
$ws1.Range("D2").Value = 46
$ws1.Range("D3").Value = 51

# Column L (Seasonality Index)
$ws1.Range("L2").Value = 0.91
$ws1.Range("L3").Value = 1.03
$ws1.Range("L4").Value = 1.19
$ws1.Range("L5").Value = 1.17
$ws1.Range("L6").Value = 1.16
$ws1.Range("L7").Value = 0.84
$ws1.Range("L8").Value = 0.99
$ws1.Range("L9").Value = 0.97
$ws1.Range("L10").Value = 1.08
$ws1.Range("L11").Value = 0.95
$ws1.Range("L12").Value = 0.88
$ws1.Range("L13").Value = 1.01
$ws1.Range("L14").Value = 1.01
$ws1.Range("L15").Value = 1.14
$ws1.Range("L16").Value = 1.13
$ws1.Range("L17").Value = 1.19

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")

# Ensure these cells keep storing text (not auto-converted to numbers)
$ws2.Range("B9:B12").NumberFormat = "@"

$ws2.Range("B9").Value = "590"
$ws2.Range("B10").Value = "320"
$ws2.Range("B11").Value = "173"
$ws2.Range("B12").Value = "51"
